$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to hold a literal text value (not auto-converted to a
    # date/number by Excel's smart-entry parsing), while leaving the cell's
    # style pointing back at the default "Normal" style once done.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Clear cells from the old layout that are no longer used anywhere in the
# new layout (old row 3 lived out at columns N:R).
$ws.Range("N3:R3").ClearContents()

# Row 1 - headers
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Rink"
$ws.Range("D1").Value = "Board Brush"
$ws.Range("E1").Value = "Wet Cut"
$ws.Range("F1").Value = "Dry Cut"
$ws.Range("G1").Value = "Edged"
$ws.Range("H1").Value = "Three Lap"
$ws.Range("I1").Value = "Flood"
$ws.Range("J1").Value = "Center Flood"
$ws.Range("K1").Value = "Dump Tank"
$ws.Range("L1").Value = "HoneyWells"
$ws.Range("M1").Value = "Room Temp/Humidity"
$ws.Range("N1").Value = "Initials"
$ws.Range("O1").Value = "Comment"

# Row 2
Set-TextValue "A2" "5/15/2017"
$ws.Range("B2").Value = "6:50PM"
$ws.Range("C2").Value = "Rink1"
$ws.Range("D2").Value = "Brush"
$ws.Range("E2").Value = "Wet"
Set-TextValue "F2" "0"
Set-TextValue "G2" "0"
Set-TextValue "H2" "0"
Set-TextValue "I2" "0"
Set-TextValue "J2" "0"

# Row 3
Set-TextValue "A3" "5/15/2017"
$ws.Range("B3").Value = "6:51PM"
$ws.Range("C3").Value = "Rink1"
$ws.Range("D3").Value = "Brush"
$ws.Range("E3").Value = "Wet"
Set-TextValue "F3" "0"
Set-TextValue "G3" "0"
Set-TextValue "H3" "0"
$ws.Range("I3").Value = "Flood"
$ws.Range("J3").Value = "Center Flood"
